# avg_dev subtract from estimate + plot
# Fill in the "Photo" counter (column D) for rows 70-102, continuing the
# existing even-number sequence (..., 142, 144, 146, 148, ..., 210).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(146,148,150,152,154,156,158,160,162,164,166,168,170,172,174,176,178,180,182,184,186,188,190,192,194,196,198,200,202,204,206,208,210)

$startRow = 70
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Scroll the sheet down to row 69 and move the active selection to H58,
# matching where the author was working when the file was last saved.
$ws.Range("H58").Select()
$excel.ActiveWindow.ScrollRow = 69
$excel.ActiveWindow.ScrollColumn = 1
